$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "applebees"
$ws.Range("B2").Value = "nonfiction"
$ws.Range("C2").Value = "1"

$ws.Range("A3").Value = "dragon ball"
$ws.Range("B3").Value = "manga"
$ws.Range("C3").Value = "1,2,3,4,5,6,7"

$ws.Range("A4").Value = "naruto"
$ws.Range("B4").Value = "manga"
$ws.Range("C4").Value = "2"

$ws.Range("A5").Value = "d"
$ws.Range("B5").Value = "man"
$ws.Range("C5").Value = "4"
